$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.78"
$ws.Range("E2").Value = "'-2.53%"
$ws.Range("G2").Value = "'18"
$ws.Range("E3").Value = "'1.38%"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'4.764"
$ws.Range("E4").Value = "'-3.64%"
$ws.Range("G4").Value = "'18"
$ws.Range("D5").Value = "'0.06308"
$ws.Range("E5").Value = "'-1.77%"
$ws.Range("G5").Value = "'18"
$ws.Range("D6").Value = "'6.921"
$ws.Range("E6").Value = "'-1.19%"
$ws.Range("G6").Value = "'18"
$ws.Range("D7").Value = "'1.349"
$ws.Range("E7").Value = "'32.99%"
$ws.Range("G7").Value = "'18"
$ws.Range("D8").Value = "'0.8772"
$ws.Range("E8").Value = "'-1.13%"
$ws.Range("G8").Value = "'18"
$ws.Range("E9").Value = "'1.24%"
$ws.Range("G9").Value = "'18"
$ws.Range("D10").Value = "'0.05011"
$ws.Range("E10").Value = "'-3.31%"
$ws.Range("G10").Value = "'18"
$ws.Range("D11").Value = "'0.07598"
$ws.Range("E11").Value = "'1.74%"
$ws.Range("G11").Value = "'18"
$ws.Range("D12").Value = "'0.02932"
$ws.Range("E12").Value = "'-5.49%"
$ws.Range("G12").Value = "'18"
$ws.Range("D13").Value = "'0.09005"
$ws.Range("E13").Value = "'-0.62%"
$ws.Range("G13").Value = "'18"
$ws.Range("D14").Value = "'0.001566"
$ws.Range("E14").Value = "'-0.23%"
$ws.Range("G14").Value = "'18"
$ws.Range("D15").Value = "'0.0006353"
$ws.Range("E15").Value = "'0.87%"
$ws.Range("G15").Value = "'18"
$ws.Range("D16").Value = "'0.005784"
$ws.Range("E16").Value = "'-3.81%"
$ws.Range("G16").Value = "'18"
$ws.Range("D17").Value = "'3.447"
$ws.Range("E17").Value = "'-1.35%"
$ws.Range("G17").Value = "'18"
$ws.Range("D18").Value = "'3.297"
$ws.Range("E18").Value = "'-1.66%"
$ws.Range("G18").Value = "'18"
$ws.Range("E19").Value = "'-0.63%"
$ws.Range("G19").Value = "'18"
$ws.Range("E20").Value = "'0.11%"
$ws.Range("G20").Value = "'18"
$ws.Range("E21").Value = "'1.05%"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'3.893"
$ws.Range("E22").Value = "'-0.29%"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'0.04405"
$ws.Range("E23").Value = "'1.08%"
$ws.Range("G23").Value = "'18"
$ws.Range("D24").Value = "'0.001169"
$ws.Range("E24").Value = "'-0.85%"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.003842"
$ws.Range("E25").Value = "'3.82%"
$ws.Range("G25").Value = "'18"
$ws.Range("D26").Value = "'0.0001198"
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("G26").Value = "'18"
$ws.Range("D27").Value = "'0.0001935"
$ws.Range("E27").Value = "'14.29%"
$ws.Range("G27").Value = "'18"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("G38").Value = "'18"
$ws.Range("G39").Value = "'18"
$ws.Range("D40").Value = "'0.04101"
$ws.Range("E40").Value = "'-1.32%"
$ws.Range("G40").Value = "'18"
$ws.Range("D41").Value = "'0.006837"
$ws.Range("E41").Value = "'2.68%"
$ws.Range("G41").Value = "'18"
$ws.Range("E42").Value = "'-0.75%"
$ws.Range("G42").Value = "'18"
$ws.Range("D43").Value = "'0.002057"
$ws.Range("E43").Value = "'-12.73%"
$ws.Range("G43").Value = "'18"
$ws.Range("D44").Value = "'0.01148"
$ws.Range("E44").Value = "'-8.94%"
$ws.Range("G44").Value = "'18"
$ws.Range("D45").Value = "'0.00005167"
$ws.Range("E45").Value = "'-1.69%"
$ws.Range("G45").Value = "'18"
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D46").Value = "'1.490"
$ws.Range("E46").Value = "'-36.76%"
$ws.Range("G46").Value = "'18"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.02299"
$ws.Range("E47").Value = "'2.20%"
$ws.Range("G47").Value = "'18"
$ws.Range("G48").Value = "'18"
$ws.Range("G49").Value = "'18"
$ws.Range("G50").Value = "'18"
$ws.Range("G51").Value = "'18"
